$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1) "Years ... : 2018 - 2020."  -- consolidate the split runs that
#    spell out "2018" and "2020" into a single run's text.
# -----------------------------------------------------------------
$r = $d.Content
$r.Find.ClearFormatting()
$r.Find.Replacement.ClearFormatting()
$r.Find.Text = ": 2018 – 2020."
$r.Find.Replacement.Text = ": 2018 – 2020."
$r.Find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null

# -----------------------------------------------------------------
# 2) "Status ... : Completed." -- merge the ": " run with "Completed."
# -----------------------------------------------------------------
$r = $d.Content
$r.Find.ClearFormatting()
$r.Find.Replacement.ClearFormatting()
$r.Find.Text = ": Completed."
$r.Find.Replacement.Text = ": Completed."
$r.Find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null

# -----------------------------------------------------------------
# 3) "Junior Data Analyst" -- merge "Junior" and " Data Analyst " runs
# -----------------------------------------------------------------
$r = $d.Content
$r.Find.ClearFormatting()
$r.Find.Replacement.ClearFormatting()
$r.Find.Text = "Junior Data Analyst "
$r.Find.Replacement.Text = "Junior Data Analyst "
$r.Find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null

# -----------------------------------------------------------------
# 4) GitHub link -- swap the username order and drop the hyperlink,
#    turning it into plain text.
# -----------------------------------------------------------------
$h = $d.Hyperlinks.Item(1)
$h.Delete()

$r = $d.Content
$r.Find.ClearFormatting()
$r.Find.Text = "https://github.com/MoretiGiven"
$r.Find.Execute() | Out-Null
$start = $r.Start
$r.Delete()
$collapsed = $d.Range($start, $start)
$collapsed.InsertAfter("https://github.com/GivenMoreti")

Write-Output "done"
